$wb = $excel.ActiveWorkbook
$sheetNames = @("展览", "全部类型")
foreach ($sheetName in $sheetNames) {
  $ws = $wb.Worksheets.Item($sheetName)

  # Ensure new rows 17-18 exist with the same formatting as row 16
  $ws.Range("A16:J16").Copy() | Out-Null
  $ws.Range("A17:J17").PasteSpecial(-4122) | Out-Null
  $ws.Range("A16:J16").Copy() | Out-Null
  $ws.Range("A18:J18").PasteSpecial(-4122) | Out-Null
  $excel.CutCopyMode = 0

  # Force text number format on text-bearing columns to avoid Excel
  # auto-converting numeric-looking / date-looking strings
  $ws.Range("B2:B18").NumberFormat = "@"
  $ws.Range("C2:C18").NumberFormat = "@"
  $ws.Range("D2:D18").NumberFormat = "@"
  $ws.Range("E2:E18").NumberFormat = "@"
  $ws.Range("G2:G18").NumberFormat = "@"
  $ws.Range("I2:I18").NumberFormat = "@"
  $ws.Range("J2:J18").NumberFormat = "@"

  # Row 2
  $ws.Cells.Item(2,1).Value = 1
  $ws.Cells.Item(2,2).Value = '2024-01-20'
  $ws.Cells.Item(2,3).Value = '合肥·第十二届次元之门动漫游戏博览会-吴磊专场'
  $ws.Cells.Item(2,4).Value = '文忠路1865号 赫拉诺言艺术中心'
  $ws.Cells.Item(2,5).Value = '2024.01.20 09:30-01.20 17:30'
  $ws.Cells.Item(2,6).Value = 422
  $ws.Cells.Item(2,7).Value = '已售罄'
  $ws.Cells.Item(2,8).Value = $false
  $ws.Cells.Item(2,9).Value = ''
  $ws.Cells.Item(2,10).Value = '//i2.hdslb.com/bfs/openplatform/202312/VBekVPuH1703840712015.jpeg'

  # Row 3
  $ws.Cells.Item(3,1).Value = 2
  $ws.Cells.Item(3,2).Value = '2024-01-27'
  $ws.Cells.Item(3,3).Value = ''
  $ws.Cells.Item(3,4).Value = '临泉路88号板桥里墨园E区1号省羽体中心 省羽体super速搏羽毛球馆'
  $ws.Cells.Item(3,5).Value = '2024.01.27 10:00-01.28 17:00'
  $ws.Cells.Item(3,6).Value = 1432
  $ws.Cells.Item(3,7).Value = '不可售'
  $ws.Cells.Item(3,8).Value = $false
  $ws.Cells.Item(3,9).Value = ''
  $ws.Cells.Item(3,10).Value = '//i1.hdslb.com/bfs/openplatform/202311/2v00jbxM1698999146733.jpeg'

  # Row 4
  $ws.Cells.Item(4,1).Value = 3
  $ws.Cells.Item(4,2).Value = '2024-01-27'
  $ws.Cells.Item(4,3).Value = '合肥·新春AG动漫游戏盛典热血plus'
  $ws.Cells.Item(4,4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
  $ws.Cells.Item(4,5).Value = '2024.01.27 10:00-01.28 17:30'
  $ws.Cells.Item(4,6).Value = 7421
  $ws.Cells.Item(4,7).Value = '65'
  $ws.Cells.Item(4,8).Value = $true
  $ws.Cells.Item(4,9).Value = ''
  $ws.Cells.Item(4,10).Value = '//i1.hdslb.com/bfs/openplatform/202312/iJ1Dnmla1702029064983.jpeg'

  # Row 5
  $ws.Cells.Item(5,1).Value = 4
  $ws.Cells.Item(5,2).Value = '2024-01-28'
  $ws.Cells.Item(5,3).Value = ''
  $ws.Cells.Item(5,4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
  $ws.Cells.Item(5,5).Value = '2024.01.28 10:00-01.28 17:00'
  $ws.Cells.Item(5,6).Value = 531
  $ws.Cells.Item(5,7).Value = '已售罄'
  $ws.Cells.Item(5,8).Value = $false
  $ws.Cells.Item(5,9).Value = ''
  $ws.Cells.Item(5,10).Value = '//i0.hdslb.com/bfs/openplatform/202312/9ClQwbVE1703668101900.jpeg'

  # Row 6
  $ws.Cells.Item(6,1).Value = 5
  $ws.Cells.Item(6,2).Value = '2024-01-28'
  $ws.Cells.Item(6,3).Value = '合肥·环形宇宙动漫游戏嘉年华—吴晛专场'
  $ws.Cells.Item(6,4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
  $ws.Cells.Item(6,5).Value = '2024.01.28 10:00-01.28 17:00'
  $ws.Cells.Item(6,6).Value = 318
  $ws.Cells.Item(6,7).Value = '258'
  $ws.Cells.Item(6,8).Value = $false
  $ws.Cells.Item(6,9).Value = ''
  $ws.Cells.Item(6,10).Value = '//i0.hdslb.com/bfs/openplatform/202312/aHzqArm61703662347629.jpeg'

  # Row 7
  $ws.Cells.Item(7,1).Value = 6
  $ws.Cells.Item(7,2).Value = '2024-01-28'
  $ws.Cells.Item(7,3).Value = '肥东· 原神&崩铁&崩坏only'
  $ws.Cells.Item(7,4).Value = '团结东路7号 巢湖宾馆'
  $ws.Cells.Item(7,5).Value = '2024.01.28 10:00-01.28 17:00'
  $ws.Cells.Item(7,6).Value = 19
  $ws.Cells.Item(7,7).Value = '55'
  $ws.Cells.Item(7,8).Value = $false
  $ws.Cells.Item(7,9).Value = ''
  $ws.Cells.Item(7,10).Value = '//i0.hdslb.com/bfs/openplatform/202401/UekMeUjQ1705462868391.jpeg'

  # Row 8
  $ws.Cells.Item(8,1).Value = 7
  $ws.Cells.Item(8,2).Value = '2024-01-29'
  $ws.Cells.Item(8,3).Value = '巢湖·原×铁×崩only'
  $ws.Cells.Item(8,4).Value = '长江东路徽商城2幢B座(祥和地铁站C口步行370米) 曼斯顿尚品酒店'
  $ws.Cells.Item(8,5).Value = '2024.01.29 10:00-01.29 17:00'
  $ws.Cells.Item(8,6).Value = 12
  $ws.Cells.Item(8,7).Value = '55'
  $ws.Cells.Item(8,8).Value = $false
  $ws.Cells.Item(8,9).Value = ''
  $ws.Cells.Item(8,10).Value = '//i0.hdslb.com/bfs/openplatform/202401/9XumHIT31705464002179.jpeg'

  # Row 9
  $ws.Cells.Item(9,1).Value = 8
  $ws.Cells.Item(9,2).Value = '2024-01-31'
  $ws.Cells.Item(9,3).Value = '巢湖·原神&崩铁&崩坏only'
  $ws.Cells.Item(9,4).Value = '仙满楼·麦肯希酒店 仙满楼·麦肯希酒店'
  $ws.Cells.Item(9,5).Value = '2024.01.31 10:00-01.31 17:00'
  $ws.Cells.Item(9,6).Value = 19
  $ws.Cells.Item(9,7).Value = '55'
  $ws.Cells.Item(9,8).Value = $false
  $ws.Cells.Item(9,9).Value = ''
  $ws.Cells.Item(9,10).Value = '//i0.hdslb.com/bfs/openplatform/202401/euD63Mlp1705479140627.jpeg'

  # Row 10
  $ws.Cells.Item(10,1).Value = 9
  $ws.Cells.Item(10,2).Value = '2024-02-03'
  $ws.Cells.Item(10,3).Value = '合肥·2024运动新春动漫庆典（全ip）'
  $ws.Cells.Item(10,4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
  $ws.Cells.Item(10,5).Value = '2024.02.03 09:30-02.04 17:00'
  $ws.Cells.Item(10,6).Value = 5417
  $ws.Cells.Item(10,7).Value = '65'
  $ws.Cells.Item(10,8).Value = $true
  $ws.Cells.Item(10,9).Value = ''
  $ws.Cells.Item(10,10).Value = '//i0.hdslb.com/bfs/openplatform/202312/tBk3WVyX1702968658234.jpeg'

  # Row 11
  $ws.Cells.Item(11,1).Value = 10
  $ws.Cells.Item(11,2).Value = '2024-02-04'
  $ws.Cells.Item(11,3).Value = '合肥·六安lovelive only'
  $ws.Cells.Item(11,4).Value = '健康东路7号 巢湖国际饭店'
  $ws.Cells.Item(11,5).Value = '2024.02.04 10:00-02.04 17:00'
  $ws.Cells.Item(11,6).Value = 8
  $ws.Cells.Item(11,7).Value = '60'
  $ws.Cells.Item(11,8).Value = $false
  $ws.Cells.Item(11,9).Value = ''
  $ws.Cells.Item(11,10).Value = '//i0.hdslb.com/bfs/openplatform/202401/wVVrdShB1705487994232.jpeg'

  # Row 12
  $ws.Cells.Item(12,1).Value = 11
  $ws.Cells.Item(12,2).Value = '2024-02-04'
  $ws.Cells.Item(12,3).Value = '合肥·梦时空SPO1动漫展'
  $ws.Cells.Item(12,4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
  $ws.Cells.Item(12,5).Value = '2024.02.04 11:30-02.04 17:00'
  $ws.Cells.Item(12,6).Value = 137
  $ws.Cells.Item(12,7).Value = '168'
  $ws.Cells.Item(12,8).Value = $false
  $ws.Cells.Item(12,9).Value = ''
  $ws.Cells.Item(12,10).Value = '//i0.hdslb.com/bfs/openplatform/202401/MSS7qIQp1704695420767.jpeg'

  # Row 13
  $ws.Cells.Item(13,1).Value = 12
  $ws.Cells.Item(13,2).Value = '2024-02-05'
  $ws.Cells.Item(13,3).Value = '合肥·国乙only新春年会版'
  $ws.Cells.Item(13,4).Value = '经开区繁华大道与莲花路交叉口 百乐门大剧院'
  $ws.Cells.Item(13,5).Value = '2024.02.05 09:00-02.05 17:00'
  $ws.Cells.Item(13,6).Value = 11
  $ws.Cells.Item(13,7).Value = '70'
  $ws.Cells.Item(13,8).Value = $false
  $ws.Cells.Item(13,9).Value = ''
  $ws.Cells.Item(13,10).Value = '//i2.hdslb.com/bfs/openplatform/202401/QkgtYncY1705656564257.jpeg'

  # Row 14
  $ws.Cells.Item(14,1).Value = 13
  $ws.Cells.Item(14,2).Value = '2024-02-13'
  $ws.Cells.Item(14,3).Value = '合肥·环形宇宙动漫游戏嘉年华'
  $ws.Cells.Item(14,4).Value = '山西路与太原路交叉口 挥动体育'
  $ws.Cells.Item(14,5).Value = '2024.02.13 09:30-02.14 16:00'
  $ws.Cells.Item(14,6).Value = 1714
  $ws.Cells.Item(14,7).Value = '39'
  $ws.Cells.Item(14,8).Value = $false
  $ws.Cells.Item(14,9).Value = ''
  $ws.Cells.Item(14,10).Value = '//i1.hdslb.com/bfs/openplatform/202401/yI94srFk1704703809648.jpeg'

  # Row 15
  $ws.Cells.Item(15,1).Value = 14
  $ws.Cells.Item(15,2).Value = '2024-02-14'
  $ws.Cells.Item(15,3).Value = '合肥·安徽马娘only'
  $ws.Cells.Item(15,4).Value = '阜阳路16号 银瑞林国际大酒店'
  $ws.Cells.Item(15,5).Value = '2024.02.14 10:00-02.14 17:00'
  $ws.Cells.Item(15,6).Value = 58
  $ws.Cells.Item(15,7).Value = '60'
  $ws.Cells.Item(15,8).Value = $false
  $ws.Cells.Item(15,9).Value = ''
  $ws.Cells.Item(15,10).Value = '//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg'

  # Row 16
  $ws.Cells.Item(16,1).Value = 15
  $ws.Cells.Item(16,2).Value = '2024-02-17'
  $ws.Cells.Item(16,3).Value = '合肥· 第二届漫画城市动漫展 -故事再次开始'
  $ws.Cells.Item(16,4).Value = '锦绣大道与清潭路交口东北角 李宁体育公园'
  $ws.Cells.Item(16,5).Value = '2024.02.17 09:00-02.17 17:00'
  $ws.Cells.Item(16,6).Value = 1086
  $ws.Cells.Item(16,7).Value = '65'
  $ws.Cells.Item(16,8).Value = $false
  $ws.Cells.Item(16,9).Value = ''
  $ws.Cells.Item(16,10).Value = '//i0.hdslb.com/bfs/openplatform/202312/vzuMc0sJ1702902061660.jpeg'

  # Row 17
  $ws.Cells.Item(17,1).Value = 16
  $ws.Cells.Item(17,2).Value = '2024-02-19'
  $ws.Cells.Item(17,3).Value = '肥西·原神&崩铁&崩坏only'
  $ws.Cells.Item(17,4).Value = '桐城路与庐江路交叉口西南80米 赤阑桥文玩大厦'
  $ws.Cells.Item(17,5).Value = '2024.02.19 09:00-02.19 17:00'
  $ws.Cells.Item(17,6).Value = 267
  $ws.Cells.Item(17,7).Value = '68'
  $ws.Cells.Item(17,8).Value = $false
  $ws.Cells.Item(17,9).Value = ''
  $ws.Cells.Item(17,10).Value = '//i1.hdslb.com/bfs/openplatform/202311/721L5pIZ1699428443216.jpeg'

  # Row 18
  $ws.Cells.Item(18,1).Value = 17
  $ws.Cells.Item(18,2).Value = '2024-04-04'
  $ws.Cells.Item(18,3).Value = '合肥·第十二届次元之门动漫游戏博览会-赵乾景专场'
  $ws.Cells.Item(18,4).Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
  $ws.Cells.Item(18,5).Value = '2024.04.04 09:00-04.05 17:00'
  $ws.Cells.Item(18,6).Value = 5476
  $ws.Cells.Item(18,7).Value = '60'
  $ws.Cells.Item(18,8).Value = $false
  $ws.Cells.Item(18,9).Value = ''
  $ws.Cells.Item(18,10).Value = '//i2.hdslb.com/bfs/openplatform/202311/244eBWip1700711342120.jpeg'

}
